# ---------------------------------------------------------------------------
# feat: Add external app
#
# - Inserts a new 'navbar.more.title' translation row right after 'navbar.log'
#   (pushes the existing channel.*/appMenu.*/... rows down by one row).
# - Renames the zh-CN text for 'channel.external' from '第三方APP' to '外部 APP'.
# - Appends two new translation rows at the bottom of the table:
#   'app.hotpot.addDish' and 'account.dialog.title'.
#
# Row/column insertion via Rows.Insert() creates brand-new cell-style (xf)
# entries in styles.xml for the shifted cells (this engine does not dedupe
# them against the existing ones), which would show up as an (unwanted) diff
# in xl/styles.xml. Instead we shift the *values* down manually (bottom-up,
# using plain .Value writes, which keep each destination cell's existing
# style intact) and only copy *formats* (PasteSpecial xlPasteFormats) for the
# brand-new rows appended past the old last row - that reuses the identical
# existing style indices too, so styles.xml stays byte-for-byte unchanged.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Snapshot of the current rows 18-62 (A/B/C), in top-to-bottom order, exactly
# as they exist before the edit. After inserting the new 'navbar.more.title'
# row at row 18, this whole block needs to move down to rows 19-63.
$oldData = @(
  @("channel.life", "生活常用", "Lifestyle"),
  @("channel.ai", "人工智能", "AI"),
  @("channel.image", "图片视频", "Media Process"),
  @("channel.developer", "编程开发", "Developer"),
  @("channel.external", "第三方APP", "External App"),
  @("channel.wip", "开发中", "WIP"),
  @("appMenu.github", "在 GitHub 上编辑", "Edit on GitHub"),
  @("appMenu.bookmark", "收藏", "Bookmark"),
  @("appMenu.frame", "框架引用", "Use in your website"),
  @("appMenu.feedback", "反馈", "Feedback"),
  @("appMenu.bookmark.undo", "移除收藏", "Remove Bookmark"),
  @("settings.language.title", "语言", "Language"),
  @("settings.language.auto", "跟随系统", "Auto"),
  @("settings.language.zh_cn", "简体中文", "简体中文"),
  @("settings.language.en_us", "English", "English"),
  @("donation.paid.title", "付费方式", "Paid Option"),
  @("feedback.send", "提交", "Submit"),
  @("feedback.hero", "我们会阅读每一条反馈", "We Read Every Feedback"),
  @("feedback.subtitle", "你可以畅所欲言", "New app request, bug report, or anything you want to tell us."),
  @("feedback.debug", "发送错误日志", "Send Error Log"),
  @("feedback.content.placeholder", "输入内容", "Write something you want to tell us"),
  @("feedback.contact.placeholder", "适合我们联系你的方式", "How can we contact you?"),
  @("general.chooseFile", "选择文件", "Choose File"),
  @("general.confirm", "确认", "Confirm"),
  @("general.save", "保存", "Save"),
  @("general.download", "下载", "Download"),
  @("qrcode.basic.title", "基本", "Basic"),
  @("qrcode.basic.placeholder", "链接或文本", "URL or Text"),
  @("qrcode.basic.type", "类型", "Type"),
  @("qrcode.basic.wifi", "WI-FI", "WI-FI"),
  @("qrcode.basic.text", "文本", "Text"),
  @("qrcode.advanced.title", "高级", "Advanced"),
  @("qrcode.advanced.icon", "图标", "Icon"),
  @("qrcode.advanced.light", "亮色", "Light Color"),
  @("qrcode.advanced.dark", "暗色", "Dark Color"),
  @("aboutPage.meta.title", "关于", "About"),
  @("app.decision.addOption", "添加选项", "Add Option"),
  @("app.decision.savePreset", "保存预设", "Save Preset"),
  @("app.decision.currentOption", "当前备选项", "Current Options"),
  @("app.roman.inputHint", "输入整数", "Input integer"),
  @("app.urlcleaner.confirmBtn", "净化", "Clean URL"),
  @("app.urlcleaner.ruleTitle", "规则", "Rules"),
  @("app.pornhub.fontSize", "字体大小", "Font Size"),
  @("app.pornhub.vertical", "竖直排列", "Vertical Layout"),
  @("app.pornhub.colorRevert", "颜色反转", "Color Revert"),
)

$firstRow = 18
$lastRow = 62
$rowCount = $lastRow - $firstRow + 1

# 1) Give the brand-new trailing rows (63-65) the same column formatting as
#    the existing data rows (A/C = style 4, B = style 5, D = style 6, E = style
#    7) plus the 20.1pt custom row height, by copying the format of the last
#    existing data row (row 62). This reuses the existing style indices rather
#    than minting new ones.
$ws.Range("A62:E62").Copy()
$newLastRow = $lastRow + 1 + 2
for ($r = $lastRow + 1; $r -le $newLastRow; $r++) {
  $ws.Range("A" + $r + ":E" + $r).PasteSpecial(-4122)
  $ws.Rows.Item($r).RowHeight = 20.1
}

# 2) Shift rows 18-62 down to 19-63, bottom row first so we never overwrite
#    data before it has been copied out (we already have it in $oldData so
#    strictly this isn't required, but keep the natural bottom-up order).
for ($i = $rowCount - 1; $i -ge 0; $i--) {
  $destRow = $firstRow + 1 + $i
  $vals = $oldData[$i]
  $ws.Cells.Item($destRow, 1).Value = $vals[0]
  $ws.Cells.Item($destRow, 2).Value = $vals[1]
  $ws.Cells.Item($destRow, 3).Value = $vals[2]
}

# 3) Write the new 'navbar.more.title' row into the now-vacated row 18 (format
#    was already correct - this row previously held 'channel.life', now moved
#    to row 19).
$ws.Cells.Item(18, 1).Value = "navbar.more.title"
$ws.Cells.Item(18, 2).Value = "更多 YGeeker 产品"
$ws.Cells.Item(18, 3).Value = "More Products from YGeeker"

# 4) Rename the zh-CN text for channel.external (now at row 23) from
#    '第三方APP' to '外部 APP'.
$ws.Cells.Item(23, 2).Value = "外部 APP"

# 5) Append the two brand-new translation rows at the bottom (64 and 65;
#    formatting for 63-65 was already applied in step 1).
$ws.Cells.Item(64, 1).Value = "app.hotpot.addDish"
$ws.Cells.Item(64, 2).Value = "从预设菜品中选择"
$ws.Cells.Item(64, 3).Value = "Add Preset Dish"

$ws.Cells.Item(65, 1).Value = "account.dialog.title"
$ws.Cells.Item(65, 2).Value = "YGeeker 账户"
$ws.Cells.Item(65, 3).Value = "YGeeker Account"

Write-Host "Applied i18n update: navbar.more.title, channel.external rename, app.hotpot.addDish, account.dialog.title"
